$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 728468337064.214
$ws.Range("C3").Value = 247582345045.7819
$ws.Range("C4").Value = 35258720994.39204
$ws.Range("C5").Value = 33256557461.85886
$ws.Range("C6").Value = 24029275225.26164
$ws.Range("C7").Value = 13494726670.04683
$ws.Range("C8").Value = 10773566732.13237
$ws.Range("C9").Value = 9012306940.781408
$ws.Range("C10").Value = 8298860851.393289
$ws.Range("C11").Value = 8107879007.03618
$ws.Range("C12").Value = 7358767498.515913
$ws.Range("C13").Value = 7212440177.912173
$ws.Range("C14").Value = 6662263144.838155
$ws.Range("C15").Value = 6051715601.171719
$ws.Range("C16").Value = 5127222348.481926
$ws.Range("C17").Value = 4765344038.656413
$ws.Range("C18").Value = 4401359373.274726
$ws.Range("C19").Value = 3778981678.160348
$ws.Range("C20").Value = 3487123027.636964
$ws.Range("C21").Value = 3293326415.479542
$ws.Range("C22").Value = 3275289772.452466
$ws.Range("C23").Value = 3054902357.104695
$ws.Range("C24").Value = 2777804281.923453
$ws.Range("C25").Value = 2721418246.808467
$ws.Range("C26").Value = 2361333413.99097
